# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D"  = 1.303
    "G"  = -0.08106060606060606
    "H"  = -0.08106060606060606
    "I"  = -0.01174242424242424
    "J"  = -0.01174242424242424
    "K"  = -0.944
    "L"  = -0.3575757575757575
    "U"  = 0.635
    "V"  = 0.1579601990049752
    "W"  = 0.9583756345177664
    "X"  = 0.3041594403114543
    "Y"  = 0.6542161942063121
    "Z"  = 0.5429864253393666
    "AA" = -0.006375976964212258
    "AB" = 0.1609555533657016
    "AC" = -0.1673315303299138
    "AD" = 6.28
    "AE" = 0
    "AF" = 6.28
    "AG" = 5.645
    "AH" = 0.6097087378640776
    "AI" = 1.04876419505678
    "AJ" = 0.5840662183135025
    "AK" = 1.054548851111526
    "AL" = 0.98
    "AM" = 0.9
    "AN" = -1570
    "AO" = -0.03163265306122449
    "AP" = -1411.25
    "AQ" = -0.03444444444444444
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
